$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This workbook (data/rename_columns.xlsx) is a lookup table with columns:
#   A = lab_id, B = column_name, C = new_name
# The edit:
#   1) Adds 3 new mapping rows for the "HLAGyn" lab (CT_N, CT_ORF1AB, Resultado)
#      right after the existing HLAGyn rows.
#   2) Renames the "DASA_2" lab_id back to "DASA" (fixing a stray "_2" suffix)
#      and fixes a copy/paste bug where the "codigo_externo_do_paciente" row
#      pointed to "codigorequisicao" instead of "requisicao".
#   3) Adds a new mapping row for "Gene ORF" -> "Ct_ORF1ab" in the DASA block.
#   4) Renames the "DB Molecular_2" lab_id back to "DB Molecular".
# ---------------------------------------------------------------------------

# Step 1: insert 3 new rows after row 21 (the last existing HLAGyn row) and
# fill them in with the new HLAGyn column mappings.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

$ws.Cells.Item(22,1).Value = "HLAGyn"
$ws.Cells.Item(22,2).Value = "CT_N"
$ws.Cells.Item(22,3).Value = "Ct_geneN"

$ws.Cells.Item(23,1).Value = "HLAGyn"
$ws.Cells.Item(23,2).Value = "CT_ORF1AB"
$ws.Cells.Item(23,3).Value = "Ct_ORF1ab"

$ws.Cells.Item(24,1).Value = "HLAGyn"
$ws.Cells.Item(24,2).Value = "Resultado"
$ws.Cells.Item(24,3).Value = "SC2_test_result"

$ws.Range("B22").Font.Name = "Calibri"
$ws.Range("B23:C24").Font.Name = "Calibri"

# Step 2: rename the "DASA_2" lab_id to "DASA" and fix the
# "codigo_externo_do_paciente" row so that it maps to "requisicao" (it
# incorrectly pointed to "codigorequisicao").
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $labId = $ws.Cells.Item($r,1).Value()
    if ($labId -eq "DASA_2") {
        $ws.Cells.Item($r,1).Value = "DASA"
        $colName = $ws.Cells.Item($r,2).Value()
        if ($colName -eq "codigo_externo_do_paciente") {
            $ws.Cells.Item($r,3).Value = "requisicao"
        }
    }
}

# Step 3: insert a new row right after the "Gene S" row of the (now renamed)
# DASA block with the new "Gene ORF" -> "Ct_ORF1ab" mapping.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $labId = $ws.Cells.Item($r,1).Value()
    $colName = $ws.Cells.Item($r,2).Value()
    if ($labId -eq "DASA" -and $colName -eq "Gene S") {
        $geneSRow = $r
    }
}
$newRow = $geneSRow + 1
$ws.Rows.Item($newRow).Insert()
$ws.Cells.Item($newRow,1).Value = "DASA"
$ws.Cells.Item($newRow,2).Value = "Gene ORF"
$ws.Cells.Item($newRow,3).Value = "Ct_ORF1ab"

$ws.Range($ws.Cells.Item($newRow,2), $ws.Cells.Item($newRow,3)).Font.Name = "Calibri"

# Step 4: rename the "DB Molecular_2" lab_id to "DB Molecular".
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $labId = $ws.Cells.Item($r,1).Value()
    if ($labId -eq "DB Molecular_2") {
        $ws.Cells.Item($r,1).Value = "DB Molecular"
    }
}

# Restore the view to a plain selection (no frozen/scrolled topLeftCell).
$ws.Range("E39").Select() | Out-Null
